$d = $word.ActiveDocument

# Remove all existing content; we'll rebuild the body paragraph-by-paragraph
# in the new order required by the edit (some text blocks are reordered,
# so an in-place split isn't enough).
$d.Content.Delete()

$sections = @(
    @{ Type = "heading"; Text = "What went well" },
    @{ Type = "normal";  Text = "I believe we did quite well as a group. There are certainly things that need to be improved, but I don’t think there could ever be a group that wouldn’t need any sort of improvement. The aspect, that I reckon, we went really well at was establishing a timeline of deadlines and distributing the work. " },
    @{ Type = "heading"; Text = "What could be improved" },
    @{ Type = "normal";  Text = "Meeting said deadlines is probably one of the things we could improve, however, we set them up with room for error and time for polishing up our work, so I reckon we stayed on track relatively well. Another thing that could be improved is establishing a common format of our work, so that we wouldn’t have to change things around in the end to keep the work consistent." },
    @{ Type = "heading"; Text = "One thing that was surprising" },
    @{ Type = "normal";  Text = "At the moment, it’s hard to properly collaborate as, a lot of the times, our timetables don’t match up. So splitting up the work evenly, working on it independently, asking for help and assisting each other, in my opinion, was the right way to go. A timeline helped us stay on track and allowed us, to my surprise, to not fall too much behind. Personally, being just out of high school, I found that my teachers pushing me to do work actually helped me. Now that no one actually cares whether I do my work or not and I have to organise myself on my own, deadlines really did help a lot. I will definitely take note of that for my future projects." },
    @{ Type = "heading"; Text = "At least one thing that you learnt about groups" },
    @{ Type = "normal";  Text = "I’ve also learnt that, in groups, a leadership role is important. In our group that role was fulfilled by Anthony and I am very thankful to him for that. He made sure we are staying on track, provided help where we needed and kept the workflow moving. He motivated and pushed me, which I believe is another big reason as to why we never fell back too far behind." }
)

$pos = 0
$count = $sections.Count
for ($i = 0; $i -lt $count; $i++) {
    $spec = $sections[$i]
    $text = $spec.Text
    $len = $text.Length

    $insertPoint = $d.Range($pos, $pos)
    $insertPoint.InsertAfter($text)

    $runRange = $d.Range($pos, $pos + $len)
    if ($spec.Type -eq "heading") {
        $runRange.Font.Italic = $true
        $runRange.Font.ItalicBi = $true
    } else {
        $runRange.Font.Italic = $false
        $runRange.Font.ItalicBi = $false
    }

    $pos = $pos + $len

    $isLast = ($i -eq ($count - 1))
    if (-not $isLast) {
        $markRange = $d.Range($pos, $pos)
        $markRange.InsertParagraphAfter()
        $pos = $pos + 1
    }
}

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
